$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 25271
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 25271
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H74").Value = 3162.2727
$ws.Range("I74").Value = 3162.2727
$ws.Range("K74").Value = 3162.2727
$ws.Range("M74").Value = -2226.2727
$ws.Range("H77").Value = 3162.2727
$ws.Range("I77").Value = 3162.2727
$ws.Range("K77").Value = 15811.3635
$ws.Range("M77").Value = -11131.3635
$ws.Range("H97").Value = 1633
$ws.Range("J97").Value = 1633
$ws.Range("L97").Value = 4899
$ws.Range("N97").Value = -5891
$ws.Range("H116").Value = 4991.375
$ws.Range("J116").Value = 4991.375
$ws.Range("L116").Value = 4991.375
$ws.Range("N116").Value = -11875.375
$ws.Range("H127").Value = 333.16666
$ws.Range("I127").Value = 333.16666
$ws.Range("K127").Value = 999.4999799999999
$ws.Range("M127").Value = 3960.50002
$ws.Range("H129").Value = 2023.8572
$ws.Range("I129").Value = 809.63635
$ws.Range("J129").Value = 3359.5
$ws.Range("K129").Value = 2428.90905
$ws.Range("L129").Value = 10078.5
$ws.Range("M129").Value = 2571.09095
$ws.Range("N129").Value = -20078.5
$ws.Range("H137").Value = 1036.2667
$ws.Range("I137").Value = 1041
$ws.Range("K137").Value = 3123
$ws.Range("M137").Value = -573
$ws.Range("H141").Value = 2517.5454
$ws.Range("I141").Value = 2517.5454
$ws.Range("K141").Value = 7552.6362
$ws.Range("M141").Value = -2372.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 3536.6667
$ws.Range("J30").Value = 3805
$ws.Range("L30").Value = 3805
$ws.Range("N30").Value = -4105
$ws.Range("H61").Value = 2888.2856
$ws.Range("I61").Value = 1357.625
$ws.Range("J61").Value = 4929.1665
$ws.Range("K61").Value = 1357.625
$ws.Range("L61").Value = 4929.1665
$ws.Range("M61").Value = -1145.625
$ws.Range("N61").Value = -5353.1665
$ws.Range("H97").Value = 1817.4445
$ws.Range("I97").Value = 1986.125
$ws.Range("K97").Value = 1986.125
$ws.Range("M97").Value = -1490.125
$ws.Range("H132").Value = 2480.3584
$ws.Range("I132").Value = 2014.6666
$ws.Range("K132").Value = 6043.9998
$ws.Range("M132").Value = -3513.9998
$ws.Range("H136").Value = 2888.2856
$ws.Range("I136").Value = 1357.625
$ws.Range("J136").Value = 4929.1665
$ws.Range("K136").Value = 4072.875
$ws.Range("L136").Value = 14787.4995
$ws.Range("M136").Value = -1522.875
$ws.Range("N136").Value = -19887.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1605.375
$ws.Range("I105").Value = 1810.3334
$ws.Range("K105").Value = 1810.3334
$ws.Range("M105").Value = -63.33339999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4669.5293
$ws.Range("I31").Value = 2447.3333
$ws.Range("J31").Value = 10002.8
$ws.Range("K31").Value = 2447.3333
$ws.Range("L31").Value = 10002.8
$ws.Range("M31").Value = -2152.3333
$ws.Range("N31").Value = -10592.8
$ws.Range("H34").Value = 4669.5293
$ws.Range("I34").Value = 2447.3333
$ws.Range("J34").Value = 10002.8
$ws.Range("K34").Value = 2447.3333
$ws.Range("L34").Value = 10002.8
$ws.Range("M34").Value = -2245.3333
$ws.Range("N34").Value = -10406.8
$ws.Range("H58").Value = 1897.25
$ws.Range("I58").Value = 1897.25
$ws.Range("K58").Value = 1897.25
$ws.Range("M58").Value = -1694.25
$ws.Range("H107").Value = 1368.56
$ws.Range("I107").Value = 1215.7693
$ws.Range("K107").Value = 1215.7693
$ws.Range("M107").Value = 704.2307000000001
$ws.Range("H134").Value = 3825.7856
$ws.Range("I134").Value = 2141.9092
$ws.Range("K134").Value = 6425.7276
$ws.Range("M134").Value = -3890.7276
$ws.Range("H136").Value = 1897.25
$ws.Range("I136").Value = 1897.25
$ws.Range("K136").Value = 5691.75
$ws.Range("M136").Value = -3141.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 114.5
$ws.Range("I7").Value = 124.14286
$ws.Range("J7").Value = 47
$ws.Range("K7").Value = 372.42858
$ws.Range("L7").Value = 141
$ws.Range("M7").Value = -260.42858
$ws.Range("N7").Value = -365
$ws.Range("H15").Value = 464.76923
$ws.Range("I15").Value = 24.375
$ws.Range("J15").Value = 1169.4
$ws.Range("K15").Value = 73.125
$ws.Range("L15").Value = 3508.2
$ws.Range("M15").Value = 66.875
$ws.Range("N15").Value = -3788.2
$ws.Range("H36").Value = 16428.428
$ws.Range("I36").Value = 16428.428
$ws.Range("K36").Value = 49285.284
$ws.Range("M36").Value = -49116.284
$ws.Range("H103").Value = 798
$ws.Range("I103").Value = 496.5
$ws.Range("J103").Value = 1099.5
$ws.Range("K103").Value = 1489.5
$ws.Range("L103").Value = 3298.5
$ws.Range("M103").Value = -610.5
$ws.Range("N103").Value = -5056.5
$ws.Range("H106").Value = 5849.381
$ws.Range("J106").Value = 5941.9
$ws.Range("L106").Value = 17825.7
$ws.Range("N106").Value = -19717.7
$ws.Range("H137").Value = 2779791.5
$ws.Range("I137").Value = 5556815.5
$ws.Range("J137").Value = 2767.3333
$ws.Range("K137").Value = 16670446.5
$ws.Range("L137").Value = 8301.999899999999
$ws.Range("M137").Value = -16665346.5
$ws.Range("N137").Value = -18501.9999
$ws.Range("H140").Value = 1589.0625
$ws.Range("I140").Value = 1278.8462
$ws.Range("K140").Value = 3836.5386
$ws.Range("M140").Value = 1343.4614

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 227908.55
$ws.Range("I2").Value = 400369.7
$ws.Range("J2").Value = 986
$ws.Range("K2").Value = 400369.7
$ws.Range("L2").Value = 986
$ws.Range("M2").Value = -400256.7
$ws.Range("N2").Value = -1212
$ws.Range("H80").Value = 3092.2666
$ws.Range("J80").Value = 3949
$ws.Range("L80").Value = 3949
$ws.Range("N80").Value = -5945
$ws.Range("H83").Value = 3092.2666
$ws.Range("J83").Value = 3949
$ws.Range("L83").Value = 19745
$ws.Range("N83").Value = -29729
$ws.Range("H106").Value = 29890
$ws.Range("J106").Value = 29890
$ws.Range("L106").Value = 29890
$ws.Range("N106").Value = -32414
$ws.Range("H132").Value = 2588.1143
$ws.Range("I132").Value = 2619.6
$ws.Range("K132").Value = 7858.799999999999
$ws.Range("M132").Value = -5328.799999999999
$ws.Range("H136").Value = 28316.186
$ws.Range("J136").Value = 28316.186
$ws.Range("L136").Value = 84948.558
$ws.Range("N136").Value = -90048.558

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9813.906999999999
$ws.Range("I7").Value = 15449.5
$ws.Range("J7").Value = 4913.391
$ws.Range("K7").Value = 15449.5
$ws.Range("L7").Value = 4913.391
$ws.Range("M7").Value = -15337.5
$ws.Range("N7").Value = -5137.391
$ws.Range("H40").Value = 7247.3335
$ws.Range("I40").Value = 6522.1304
$ws.Range("K40").Value = 6522.1304
$ws.Range("M40").Value = -6386.1304
$ws.Range("H55").Value = 2511.8064
$ws.Range("I55").Value = 2618.0625
$ws.Range("J55").Value = 2398.4666
$ws.Range("K55").Value = 2618.0625
$ws.Range("L55").Value = 2398.4666
$ws.Range("M55").Value = -2445.0625
$ws.Range("N55").Value = -2744.4666
$ws.Range("H74").Value = 21738.666
$ws.Range("J74").Value = 22608
$ws.Range("L74").Value = 22608
$ws.Range("N74").Value = -24604
$ws.Range("H77").Value = 21738.666
$ws.Range("J77").Value = 22608
$ws.Range("L77").Value = 67824
$ws.Range("N77").Value = -77808
$ws.Range("H81").Value = 25181
$ws.Range("J81").Value = 25181
$ws.Range("L81").Value = 25181
$ws.Range("N81").Value = -27177
$ws.Range("H84").Value = 25181
$ws.Range("J84").Value = 25181
$ws.Range("L84").Value = 75543
$ws.Range("N84").Value = -85527
$ws.Range("H93").Value = 16033.167
$ws.Range("I93").Value = 1996.7368
$ws.Range("J93").Value = 69371.60000000001
$ws.Range("K93").Value = 1996.7368
$ws.Range("L93").Value = 69371.60000000001
$ws.Range("M93").Value = -748.7367999999999
$ws.Range("N93").Value = -71867.60000000001
$ws.Range("H122").Value = 186331.81
$ws.Range("I122").Value = 502501
$ws.Range("K122").Value = 1507503
$ws.Range("M122").Value = -1505053
$ws.Range("H126").Value = 9813.906999999999
$ws.Range("I126").Value = 15449.5
$ws.Range("J126").Value = 4913.391
$ws.Range("K126").Value = 46348.5
$ws.Range("L126").Value = 14740.173
$ws.Range("M126").Value = -43878.5
$ws.Range("N126").Value = -19680.173
$ws.Range("H136").Value = 4099.6
$ws.Range("I136").Value = 1000.3333
$ws.Range("J136").Value = 5427.857
$ws.Range("K136").Value = 3000.9999
$ws.Range("L136").Value = 16283.571
$ws.Range("M136").Value = -450.9998999999998
$ws.Range("N136").Value = -21383.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 13500
$ws.Range("J76").Value = 13500
$ws.Range("L76").Value = 13500
$ws.Range("N76").Value = -14130
$ws.Range("H79").Value = 13500
$ws.Range("J79").Value = 13500
$ws.Range("L79").Value = 13500
$ws.Range("N79").Value = -15684
